$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.076.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.812.20'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.44%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.51'
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4625'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.90%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3758'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.71%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07402'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8631'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.57'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.66%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.814.74'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.646'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.384'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07070'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.90'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.73%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008724'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("E19").Value = '  +0.11%  '

$ws.Range("E20").Value = '  +0.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.100.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.319'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.98%  '

$ws.Range("E23").Value = '  +0.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.042.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.96%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.915'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.193'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.48'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.259'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.97%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.65%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08919'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7728'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.90%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.173'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.520'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.905'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.55%  '

$ws.Range("E36").Value = '  +0.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.129'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01955'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05231'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.928'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.238'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.35%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5286'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.58%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.352'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +14.77%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1677'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.598'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5024'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.45%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.673'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.38%  '

$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.001'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.11%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06332'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.01%  '
